# Add a new cage (RackID 2 / CageID 5) that has no session data yet —
# i.e. a "cage with no sessions" — as a new row at the bottom of the
# "Data" sheet, and move the active selection down past it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New row 7: RackID=2, CageID=5, is Dummy=FALSE (no session columns filled in).
$ws.Cells.Item(7, 1).Value = 2
$ws.Cells.Item(7, 2).Value = 5

# Match the boolean display style ("TRUE"/"FALSE") already used by the
# other rows in the "is Dummy" column before writing the boolean value.
$ws.Range("C7").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Cells.Item(7, 3).Value = $false

# Move / record the active selection on the Data sheet.
[void]$ws.Activate()
[void]$ws.Range("F19").Select()
